# Weekly update: insert a new week's data (fecha serial 44620) at rows 16-17,
# shifting all the existing historical rows down by two and extending the
# sheet from A1:R71 to A1:R73.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 16 (the oldest-dated
# block after the most-recent-dated rows 2-15). Excel automatically shifts
# every row from the old 16..71 down to 18..73 and grows the used range.
$ws.Rows("16:17").Insert()

# Build the two new data rows (columns A..R).
$newRows = New-Object 'object[,]' 2,18

# Row 16 - Primera
$newRows[0,0]  = 9
$newRows[0,1]  = "Vega Central Mapocho de Santiago"
$newRows[0,2]  = "Metropolitana"
$newRows[0,3]  = 44620
$newRows[0,4]  = 13
$newRows[0,5]  = 100114002
$newRows[0,6]  = "Camote"
$newRows[0,7]  = "Sin especificar"
$newRows[0,8]  = "Primera"
$newRows[0,9]  = 1150
$newRows[0,10] = 11000
$newRows[0,11] = 12000
$newRows[0,12] = 11500
$newRows[0,13] = "`$/malla 18 kilos"
$newRows[0,14] = "Perú"
$newRows[0,15] = 639
$newRows[0,16] = 18
$newRows[0,17] = "Hortaliza"

# Row 17 - Segunda
$newRows[1,0]  = 9
$newRows[1,1]  = "Vega Central Mapocho de Santiago"
$newRows[1,2]  = "Metropolitana"
$newRows[1,3]  = 44620
$newRows[1,4]  = 13
$newRows[1,5]  = 100114002
$newRows[1,6]  = "Camote"
$newRows[1,7]  = "Sin especificar"
$newRows[1,8]  = "Segunda"
$newRows[1,9]  = 520
$newRows[1,10] = 9000
$newRows[1,11] = 9000
$newRows[1,12] = 9000
$newRows[1,13] = "`$/malla 18 kilos"
$newRows[1,14] = "Perú"
$newRows[1,15] = 500
$newRows[1,16] = 18
$newRows[1,17] = "Hortaliza"

$ws.Range("A16:R17").Value = $newRows
